$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.245.06"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.28"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6708"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07420"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2963"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.79"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07723"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.025"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6783"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.756.16"
$ws.Range("E14").Value = "  -4.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.44"
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.179"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008313"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.006.93"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.89"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.180"
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9996"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.95"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.700"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1398"
$ws.Range("E26").Value = "  -4.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.03"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.087"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.199"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05346"
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7606"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.877"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.144"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.676"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.335.14"
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01807"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.730"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9228"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("E41").Value = "  +2.92%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("B43").Value = "XinFinNetwork"
$ws.Range("C43").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08233"
$ws.Range("E43").Value = "  +16.90%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.34"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000125"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5160"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.767"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.914.85"
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.275"
$ws.Range("E50").Value = "  -4.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05926"
$ws.Range("E51").Value = "  +0.05%  "
